$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 (table row 1)
$t.Cell(1, 1).Range.Text = "99÷8="
$t.Cell(1, 2).Range.Text = "96÷9="
$t.Cell(1, 3).Range.Text = "27÷3="
$t.Cell(1, 4).Range.Text = "80÷8="
$t.Cell(1, 5).Range.Text = "63÷5="

# Row 5 (table row 5)
$t.Cell(5, 1).Range.Text = "85÷4="
$t.Cell(5, 2).Range.Text = "85÷9="
$t.Cell(5, 3).Range.Text = "17÷5="
$t.Cell(5, 4).Range.Text = "88÷2="
$t.Cell(5, 5).Range.Text = "36÷4="

# Row 9 (table row 9)
$t.Cell(9, 1).Range.Text = "94÷6="
$t.Cell(9, 2).Range.Text = "38÷9="
$t.Cell(9, 3).Range.Text = "80÷8="
$t.Cell(9, 4).Range.Text = "87÷8="
$t.Cell(9, 5).Range.Text = "84÷5="

# Row 13 (table row 13)
$t.Cell(13, 1).Range.Text = "28÷5="
$t.Cell(13, 2).Range.Text = "24÷7="
$t.Cell(13, 3).Range.Text = "61÷3="
$t.Cell(13, 4).Range.Text = "86÷3="
$t.Cell(13, 5).Range.Text = "74÷8="

# Row 17 (table row 17)
$t.Cell(17, 1).Range.Text = "43÷7="
$t.Cell(17, 2).Range.Text = "85÷7="
$t.Cell(17, 3).Range.Text = "32÷2="
$t.Cell(17, 4).Range.Text = "42÷6="
$t.Cell(17, 5).Range.Text = "59÷3="
